$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns for the rows being updated so that
# numeric-looking strings (e.g. "0.999", "32.76") are preserved verbatim as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E50").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '34.667.63'
$ws.Range("E2").Value = '  +1.57%  '

# Row 3
$ws.Range("D3").Value = '1.807.56'
$ws.Range("E3").Value = '  +0.90%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("E5").Value = '  -0.75%  '

# Row 6
$ws.Range("E6").Value = '  -0.44%  '

# Row 7
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").Value = '32.76'
$ws.Range("E8").Value = '  +4.63%  '

# Row 9
$ws.Range("E9").Value = '  +3.09%  '

# Row 10
$ws.Range("D10").Value = '0.0711'
$ws.Range("E10").Value = '  +7.54%  '

# Row 11
$ws.Range("D11").Value = '0.0930'
$ws.Range("E11").Value = '  +0.21%  '

# Row 12
$ws.Range("D12").Value = '2.067.65'
$ws.Range("E12").Value = '  +1.03%  '

# Row 13
$ws.Range("D13").Value = '11.13'
$ws.Range("E13").Value = '  -2.67%  '

# Row 14
$ws.Range("D14").Value = '1.813.14'
$ws.Range("E14").Value = '  +1.38%  '

# Row 15
$ws.Range("D15").Value = '0.645'
$ws.Range("E15").Value = '  +1.29%  '

# Row 16
$ws.Range("D16").Value = '34.683.22'
$ws.Range("E16").Value = '  +1.69%  '

# Row 17
$ws.Range("E17").Value = '  +2.38%  '

# Row 18
$ws.Range("D18").Value = '69.55'
$ws.Range("E18").Value = '  -0.07%  '

# Row 19
$ws.Range("D19").Value = '254.19'
$ws.Range("E19").Value = '  +0.25%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0803'
$ws.Range("E20").Value = '  +7.99%  '

# Row 21
$ws.Range("D21").Value = '11.01'
$ws.Range("E21").Value = '  +5.11%  '

# Row 22
$ws.Range("E22").Value = '  -0.03%  '

# Row 23
$ws.Range("E23").Value = '  -0.84%  '

# Row 24
$ws.Range("D24").Value = '2.18'
$ws.Range("E24").Value = '  +1.63%  '

# Row 25
$ws.Range("D25").Value = '161.68'
$ws.Range("E25").Value = '  +3.11%  '

# Row 26
$ws.Range("D26").Value = '16.49'
$ws.Range("E26").Value = '  -0.81%  '

# Row 27
$ws.Range("E27").Value = '  +1.61%  '

# Row 28
$ws.Range("E28").Value = '  -0.05%  '

# Row 29
$ws.Range("D29").Value = '677.51'
$ws.Range("E29").Value = '  +1,190.76%  '

# Row 30
$ws.Range("E30").Value = '  -0.02%  '

# Row 31
$ws.Range("E31").Value = '  +3.11%  '

# Row 32
$ws.Range("D32").Value = '3.80'
$ws.Range("E32").Value = '  -0.71%  '

# Row 33
$ws.Range("E33").Value = '  +0.27%  '

# Row 34
$ws.Range("D34").Value = '3.65'
$ws.Range("E34").Value = '  +0.69%  '

# Row 35
$ws.Range("E35").Value = '  +2.45%  '

# Row 36
$ws.Range("D36").Value = '1.438.80'
$ws.Range("E36").Value = '  -1.00%  '

# Row 37
$ws.Range("E37").Value = '  -0.33%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '0.645'
$ws.Range("E38").Value = '  +1.71%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.0193'
$ws.Range("E39").Value = '  +3.05%  '

# Row 40
$ws.Range("D40").Value = '85.16'
$ws.Range("E40").Value = '  +2.02%  '

# Row 41
$ws.Range("D41").Value = '0.956'
$ws.Range("E41").Value = '  +5.94%  '

# Row 42
$ws.Range("E42").Value = '  -1.27%  '

# Row 43
$ws.Range("E43").Value = '  -0.06%  '

# Row 44
$ws.Range("E44").Value = '  +3.22%  '

# Row 45
$ws.Range("E45").Value = '  +5.07%  '

# Row 46
$ws.Range("E46").Value = '  -1.03%  '

# Row 47
$ws.Range("D47").Value = '0.0495'
$ws.Range("E47").Value = '  -3.17%  '

# Row 48
$ws.Range("D48").Value = '1.962.91'
$ws.Range("E48").Value = '  +0.87%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '12.29'
$ws.Range("E49").Value = '  +3.06%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '106.28'
$ws.Range("E50").Value = '  +8.71%  '
